$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing "Stock category" column (C) so it
# becomes column E; the two new columns will hold the French / English common
# names for each stock.
$ws.Columns("C:D").Insert()

# Header row
$ws.Cells.Item(1, 3).Value = "Name_FR"
$ws.Cells.Item(1, 4).Value = "Name_EN"

# Common-name lookup, keyed by row number. Rows without a known common name
# (e.g. mixed-species / multi-species stocks) are intentionally left blank.
$names = @(
    @{Row=2; FR="Ange de mer commun"; EN="Angelshark"},
    @{Row=3; FR="Requin-pèlerin"; EN="Basking shark"},
    @{Row=5; FR="Aiguillat commun"; EN="Picked dogfish"},
    @{Row=6; FR="Requin hâ"; EN="Tope shark"},
    @{Row=7; FR="Squale-chagrin de l'Atlantique"; EN="Leafscale gulper shark"},
    @{Row=8; FR="Requin-taupe commun"; EN="Porbeagle"},
    @{Row=13; FR="Raie blanche"; EN="White skate"},
    @{Row=14; FR="Pocheteau gris"; EN="Blue skate"},
    @{Row=15; FR="Pocheteau gris"; EN="Blue skate"},
    @{Row=16; FR="Pocheteau gris"; EN="Blue skate"},
    @{Row=17; FR="Raie bouclée"; EN="Thornback ray"},
    @{Row=18; FR="Raie bouclée"; EN="Thornback ray"},
    @{Row=19; FR="Raie bouclée"; EN="Thornback ray"},
    @{Row=20; FR="Raie bouclée"; EN="Thornback ray"},
    @{Row=21; FR="Raie bouclée"; EN="Thornback ray"},
    @{Row=22; FR="Raie bouclée"; EN="Thornback ray"},
    @{Row=23; FR="Raie mêlée"; EN="Small-eyed ray"},
    @{Row=24; FR="Raie mêlée"; EN="Small-eyed ray"},
    @{Row=25; FR="Raie chardon"; EN="Shagreen ray"},
    @{Row=26; FR="Raie lisse"; EN="Blonde ray"},
    @{Row=27; FR="Raie lisse"; EN="Blonde ray"},
    @{Row=28; FR="Raie lisse"; EN="Blonde ray"},
    @{Row=29; FR="Raie lisse"; EN="Blonde ray"},
    @{Row=30; FR="Raie lisse"; EN="Blonde ray"},
    @{Row=31; FR="Raie circulaire"; EN="Sandy ray"},
    @{Row=32; FR="Raie douce"; EN="Spotted ray"},
    @{Row=33; FR="Raie douce"; EN="Spotted ray"},
    @{Row=34; FR="Raie douce"; EN="Spotted ray"},
    @{Row=35; FR="Raie douce"; EN="Spotted ray"},
    @{Row=36; FR="Raie douce"; EN="Spotted ray"},
    @{Row=37; FR="Raie fleurie"; EN="Cuckoo ray"},
    @{Row=38; FR="Raie fleurie"; EN="Cuckoo ray"},
    @{Row=39; FR="Raie fleurie"; EN="Cuckoo ray"},
    @{Row=40; FR="Raie fleurie"; EN="Cuckoo ray"},
    @{Row=41; FR="Raie radiée épineuse"; EN="Starry ray"},
    @{Row=42; FR="Raie brunette"; EN="Undulate ray"},
    @{Row=44; FR="Raie brunette"; EN="Undulate ray"},
    @{Row=45; FR="Raie brunette"; EN="Undulate ray"},
    @{Row=46; FR="Raie brunette"; EN="Undulate ray"},
    @{Row=47; FR="Raie brunette"; EN="Undulate ray"},
    @{Row=48; FR="Squale-liche"; EN="Kitefin shark"},
    @{Row=49; FR="Émissole tachetée"; EN="Starry smooth-hound"},
    @{Row=50; FR="Chien espagnol"; EN="Blackmouth catshark"},
    @{Row=51; FR="Chien espagnol"; EN="Blackmouth catshark"},
    @{Row=52; FR="Petite roussette"; EN="Lesser spotted dogfish"},
    @{Row=53; FR="Petite roussette"; EN="Lesser spotted dogfish"},
    @{Row=54; FR="Petite roussette"; EN="Lesser spotted dogfish"},
    @{Row=55; FR="Petite roussette"; EN="Lesser spotted dogfish"},
    @{Row=56; FR="Grande roussette"; EN="Nursehound"}
)

foreach ($entry in $names) {
    $ws.Cells.Item($entry.Row, 3).Value = $entry.FR
    $ws.Cells.Item($entry.Row, 4).Value = $entry.EN
}
